$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.378.73"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "1.885.32"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'0.698"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'246.73"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'43.28"
$ws.Range("E8").Value = "  +5.63%  "
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").Value = "'0.0749"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").Value = "'0.0980"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "'13.57"
$ws.Range("E12").Value = "  +6.21%  "
$ws.Range("D13").Value = "'0.775"
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").Value = "2.159.30"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'4.96"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "1.876.31"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "35.352.89"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "'73.79"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "0.0₃0830"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'245.42"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "'12.87"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  +6.89%  "
$ws.Range("D23").Value = "'2.63"
$ws.Range("E23").Value = "  +10.33%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "'164.86"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'8.69"
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Value = "'0.0597"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'1.47"
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("E36").Value = "  +3.88%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'0.0741"
$ws.Range("E38").Value = "  +11.84%  "
$ws.Range("D39").Value = "'17.31"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("D41").Value = "'97.74"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").Value = "1.313.84"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "'0.0806"
$ws.Range("E45").Value = "  +6.21%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'12.13"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'6.35"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").Value = "'42.65"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "2.063.39"
$ws.Range("E51").Value = "  +0.07%  "
